$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.328.55"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.856.46"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").Value = "313.72"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("D7").Value = "0.4609"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").Value = "0.3709"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.07321"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "0.8806"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "0.07805"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "1.840.72"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "5.389"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "6.543"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "91.82"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "0.000009059"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "27.356.07"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "5.130"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "10.52"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "2.148.05"
$ws.Range("E24").Value = "  +7.44%  "
$ws.Range("D25").Value = "1.939"
$ws.Range("E25").Value = "  +5.73%  "
$ws.Range("D26").Value = "152.03"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "18.40"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").Value = "5.103"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "116.16"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "0.08863"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "0.7716"
$ws.Range("E32").Value = "  +6.22%  "
$ws.Range("D33").Value = "3.037"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "1.173"
$ws.Range("E34").Value = "  +3.34%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").Value = "2.631"
$ws.Range("E36").Value = "  +5.78%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "0.01959"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "0.05232"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "2.951"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "7.048"
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("D42").Value = "0.5139"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "8.396"
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("D45").Value = "0.4831"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "10.34"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").Value = "103.30"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "1.651"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "0.06220"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "65.97"
$ws.Range("E51").Value = "  +2.13%  "
